$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    # Set slide background to solid white fill.
    $s.Background.Fill.Solid()
    $s.Background.Fill.ForeColor.RGB = 16777215

    # Style the title run: black text, Arial font.
    $title = $s.Shapes.Item(1)
    $tr = $title.TextFrame.TextRange
    $tr.Font.Color.RGB = 0
    $tr.Font.Name = "Arial"
}
